# Applies: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Reorders/refreshes the "Periodo Mora" detail rows (16-39) for the two
# trabajadores (ANDI JOSE BANQUEZ JULIO / SANTIAGO FERNANDO SANCHEZ GOMEZ)
# so that each worker's periods are newest-first, and refreshes the
# associated Valor Mora / Salario Basico figures to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Worker 1: ANDI JOSE BANQUEZ JULIO (CC 1047474273) -> rows 16-24
# Periods newest (1801) to oldest (1705), mora=29509, salario=737717
$andiDoc = "CC"
$andiId = "1047474273"
$andiName = "ANDI JOSE BANQUEZ JULIO"
$andiPeriods = @("1801","1712","1711","1710","1709","1708","1707","1706","1705")
$andiMora = 29509
$andiSalario = 737717

$row = 16
foreach ($p in $andiPeriods) {
    $ws.Cells.Item($row, 2).Value = $andiDoc
    $ws.Cells.Item($row, 3).Value = $andiId
    $ws.Cells.Item($row, 4).Value = $andiName
    $ws.Cells.Item($row, 5).Value = $p
    $ws.Cells.Item($row, 6).Value = $andiMora
    $ws.Cells.Item($row, 7).Value = $andiSalario
    $row = $row + 1
}

# Worker 2: SANTIAGO FERNANDO SANCHEZ GOMEZ (CC 79328825) -> rows 25-39
# Periods newest (2003) to oldest (1710), mora=40000 (except 2003 -> 38666),
# salario=1000000
$santiagoDoc = "CC"
$santiagoId = "79328825"
$santiagoName = "SANTIAGO FERNANDO SANCHEZ GOMEZ"
$santiagoPeriods = @("2003","2002","2001","1912","1911","1910","1909","1908","1907","1906","1905","1801","1712","1711","1710")
$santiagoSalario = 1000000

$row = 25
foreach ($p in $santiagoPeriods) {
    if ($p -eq "2003") {
        $mora = 38666
    } else {
        $mora = 40000
    }
    $ws.Cells.Item($row, 2).Value = $santiagoDoc
    $ws.Cells.Item($row, 3).Value = $santiagoId
    $ws.Cells.Item($row, 4).Value = $santiagoName
    $ws.Cells.Item($row, 5).Value = $p
    $ws.Cells.Item($row, 6).Value = $mora
    $ws.Cells.Item($row, 7).Value = $santiagoSalario
    $row = $row + 1
}
